$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C24:D44").NumberFormat = "@"
$data = New-Object "object[,]" 21,6
$data[0,0] = "Fahey Herman"
$data[0,1] = "Alfonzo"
$data[0,2] = "8652538063"
$data[0,3] = "5572237379"
$data[0,4] = "FaheyHerman@yopmail.com"
$data[0,5] = "AD4493"
$data[1,0] = "O'Keefe Hammes"
$data[1,1] = "Eddie"
$data[1,2] = "8724590647"
$data[1,3] = "3990586038"
$data[1,4] = "O'KeefeHammes@yopmail.com"
$data[1,5] = "AD3511"
$data[2,0] = "McLaughlin Pouros"
$data[2,1] = "Quinton"
$data[2,2] = "4550718771"
$data[2,3] = "5874425139"
$data[2,4] = "McLaughlinPouros@yopmail.com"
$data[2,5] = "AD1711"
$data[3,0] = "Hayes Hansen"
$data[3,1] = "Sidney"
$data[3,2] = "3885002395"
$data[3,3] = "4337279910"
$data[3,4] = "HayesHansen@yopmail.com"
$data[3,5] = "AD3445"
$data[4,0] = "Connelly Ritchie"
$data[4,1] = "Cole"
$data[4,2] = "7120221180"
$data[4,3] = "5641005536"
$data[4,4] = "ConnellyRitchie@yopmail.com"
$data[4,5] = "AD5555"
$data[5,0] = "Thompson Harber"
$data[5,1] = "Jordan"
$data[5,2] = "9029748010"
$data[5,3] = "3728357261"
$data[5,4] = "ThompsonHarber@yopmail.com"
$data[5,5] = "AD5779"
$data[6,0] = "Weimann Kunde"
$data[6,1] = "Fausto"
$data[6,2] = "3827946282"
$data[6,3] = "9856075299"
$data[6,4] = "WeimannKunde@yopmail.com"
$data[6,5] = "AD5793"
$data[7,0] = "Walsh Graham"
$data[7,1] = "Derrick"
$data[7,2] = "3189179555"
$data[7,3] = "6971637657"
$data[7,4] = "WalshGraham@yopmail.com"
$data[7,5] = "AD7849"
$data[8,0] = "Lakin Klein"
$data[8,1] = "Gabriel"
$data[8,2] = "9769420742"
$data[8,3] = "5980466489"
$data[8,4] = "LakinKlein@yopmail.com"
$data[8,5] = "AD4296"
$data[9,0] = "Parker Wilkinson"
$data[9,1] = "Denita"
$data[9,2] = "6860817531"
$data[9,3] = "8425177727"
$data[9,4] = "ParkerWilkinson@yopmail.com"
$data[9,5] = "AD4090"
$data[10,0] = "Schiller Haag"
$data[10,1] = "Rosena"
$data[10,2] = "6750158246"
$data[10,3] = "7985191572"
$data[10,4] = "SchillerHaag@yopmail.com"
$data[10,5] = "AD2690"
$data[11,0] = "Pacocha Predovic"
$data[11,1] = "Son"
$data[11,2] = "5880252142"
$data[11,3] = "5558094418"
$data[11,4] = "PacochaPredovic@yopmail.com"
$data[11,5] = "AD8359"
$data[12,0] = "Reynolds Dooley"
$data[12,1] = "Kenneth"
$data[12,2] = "7715855000"
$data[12,3] = "9011948470"
$data[12,4] = "ReynoldsDooley@yopmail.com"
$data[12,5] = "AD4765"
$data[13,0] = "Welch Beatty"
$data[13,1] = "Herschel"
$data[13,2] = "7084097873"
$data[13,3] = "4532690656"
$data[13,4] = "WelchBeatty@yopmail.com"
$data[13,5] = "AD4868"
$data[14,0] = "Anderson Dach"
$data[14,1] = "Torie"
$data[14,2] = "7125453977"
$data[14,3] = "4459052272"
$data[14,4] = "AndersonDach@yopmail.com"
$data[14,5] = "AD1517"
$data[15,0] = "Ernser Prosacco"
$data[15,1] = "Rich"
$data[15,2] = "4531675124"
$data[15,3] = "4234642624"
$data[15,4] = "ErnserProsacco@yopmail.com"
$data[15,5] = "AD8397"
$data[16,0] = "Lowe Koelpin"
$data[16,1] = "Lazaro"
$data[16,2] = "3549116738"
$data[16,3] = "7448602999"
$data[16,4] = "LoweKoelpin@yopmail.com"
$data[16,5] = "AD4009"
$data[17,0] = "Luettgen Johns"
$data[17,1] = "Traci"
$data[17,2] = "3172414030"
$data[17,3] = "9039385988"
$data[17,4] = "LuettgenJohns@yopmail.com"
$data[17,5] = "AD3406"
$data[18,0] = "Schiller Harvey"
$data[18,1] = "Jamie"
$data[18,2] = "5377584625"
$data[18,3] = "4297897836"
$data[18,4] = "SchillerHarvey@yopmail.com"
$data[18,5] = "AD8192"
$data[19,0] = "Krajcik Schneider"
$data[19,1] = "Onie"
$data[19,2] = "3634201583"
$data[19,3] = "8761983187"
$data[19,4] = "KrajcikSchneider@yopmail.com"
$data[19,5] = "AD6229"
$data[20,0] = "Waelchi Kovacek"
$data[20,1] = "Layne"
$data[20,2] = "4621243652"
$data[20,3] = "8853880683"
$data[20,4] = "WaelchiKovacek@yopmail.com"
$data[20,5] = "AD2316"
$ws.Range("A24:F44").Value = $data
